# Updated ablaufplan and systemtest.
# Replaces the fingerprint-scan based login/voting confirmation flow with an
# RFID-tag based flow, and refreshes the sheet view position/zoom/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D4").Value  = 'Öffnen der Software -> Login-Daten eingeben -> RFID-Tag einlesen'
$ws.Range("D5").Value  = 'Öffnen der Software -> Falsche Login-Daten eingeben -> RFID-Tag einlesen'
$ws.Range("E5").Value  = 'Login schlägt fehl und Fehlermeldung wird angezeigt. "Falsche Login Daten"'

$ws.Range("D22").Value = 'Bundestagswahl aus Liste der verfügbaren Wahlen auswählen -> Erststimme abgeben -> Zweitstimme abgeben -> Button "Auswahl bestätigen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'
$ws.Range("D23").Value = 'Bundestagswahl aus Liste der verfügbaren Wahlen auswählen -> Button "Stimme ungültig machen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'

$ws.Range("D25").Value = 'Europawahl aus Liste der verfügbaren Wahlen auswählen -> Stimme abgeben -> Button "Auswahl bestätigen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'
$ws.Range("D26").Value = 'Europawahl aus Liste der verfügbaren Wahlen auswählen -> Button "Stimme ungültig machen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'

$ws.Range("D28").Value = 'Bürgerentscheid aus Liste der verfügbaren Wahlen auswählen -> Stimme abgeben -> Button "Auswahl bestätigen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'
$ws.Range("D29").Value = 'Bürgerentscheid aus Liste der verfügbaren Wahlen auswählen -> Button "Stimme ungültig machen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'

$ws.Range("D31").Value = 'Landratswahl aus Liste der verfügbaren Wahlen auswählen -> Stimme abgeben -> Button "Auswahl bestätigen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'
$ws.Range("D32").Value = 'Landratswahl aus Liste der verfügbaren Wahlen auswählen -> Button "Stimme ungültig machen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'

$ws.Range("D34").Value = 'Bürgermeisterwahl aus Liste der verfügbaren Wahlen auswählen -> Stimme mehrfach oder gar nicht abgeben -> Button "Auswahl bestätigen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'
$ws.Range("D35").Value = 'Bürgermeisterwahl aus Liste der verfügbaren Wahlen auswählen -> Button "Stimme ungültig machen" anklicken -> RFID-Tag einscannen -> "Bestätigen" anklicken'

# Update the active sheet view: scroll position, zoom level and current selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 96
$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D35").Select()
